$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the match data (columns F:V) between rows 5 and 6.
#    Columns A:E (index / country / tournament / season / match date) stay put.
# ---------------------------------------------------------------------------
$row5 = @("MC Alger", 4, "Ben Aknoun", 0, 1.4, "16/09/2023 03:43", 1.3, "16/09/2023 10:40", 4.19, "16/09/2023 03:43", 4.81, "16/09/2023 16:47", 8.529999999999999, "16/09/2023 03:43", 12.64, "16/09/2023 16:47", "https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-es-ben-aknoun/WjyqCu9h/")
$row6 = @("Magra", 0, "Kabylie", 1, 2.84, "15/09/2023 13:42", 3.8, "16/09/2023 16:12", 2.63, "15/09/2023 13:42", 2.84, "16/09/2023 15:03", 2.72, "15/09/2023 13:42", 2.26, "16/09/2023 16:12", "https://www.betexplorer.com/football/algeria/ligue-1/magra-kabylie/YFXa8c8H/")

$ws.Range("F5").Value = $row6[0]
$ws.Range("G5").Value = $row6[1]
$ws.Range("H5").Value = $row6[2]
$ws.Range("I5").Value = $row6[3]
$ws.Range("J5").Value = $row6[4]
$ws.Range("K5").Value = $row6[5]
$ws.Range("L5").Value = $row6[6]
$ws.Range("M5").Value = $row6[7]
$ws.Range("N5").Value = $row6[8]
$ws.Range("O5").Value = $row6[9]
$ws.Range("P5").Value = $row6[10]
$ws.Range("Q5").Value = $row6[11]
$ws.Range("R5").Value = $row6[12]
$ws.Range("S5").Value = $row6[13]
$ws.Range("T5").Value = $row6[14]
$ws.Range("U5").Value = $row6[15]
$ws.Range("V5").Value = $row6[16]

$ws.Range("F6").Value = $row5[0]
$ws.Range("G6").Value = $row5[1]
$ws.Range("H6").Value = $row5[2]
$ws.Range("I6").Value = $row5[3]
$ws.Range("J6").Value = $row5[4]
$ws.Range("K6").Value = $row5[5]
$ws.Range("L6").Value = $row5[6]
$ws.Range("M6").Value = $row5[7]
$ws.Range("N6").Value = $row5[8]
$ws.Range("O6").Value = $row5[9]
$ws.Range("P6").Value = $row5[10]
$ws.Range("Q6").Value = $row5[11]
$ws.Range("R6").Value = $row5[12]
$ws.Range("S6").Value = $row5[13]
$ws.Range("T6").Value = $row5[14]
$ws.Range("U6").Value = $row5[15]
$ws.Range("V6").Value = $row5[16]

# ---------------------------------------------------------------------------
# 2) Swap the match data (columns F:V) between rows 17 and 18.
# ---------------------------------------------------------------------------
$row17 = @("Khenchela", 2, "Kabylie", 1, 2.63, "28/09/2023 04:12", 2.05, "29/09/2023 16:41", 2.62, "28/09/2023 04:12", 2.75, "29/09/2023 16:41", 3.02, "28/09/2023 04:12", 4.88, "29/09/2023 16:27", "https://www.betexplorer.com/football/algeria/ligue-1/khenchela-kabylie/pUZYGLcr/")
$row18 = @("US Souf", 0, "Oran", 0, 2.49, "28/09/2023 19:27", 2.14, "29/09/2023 13:29", 2.88, "28/09/2023 19:27", 2.74, "29/09/2023 14:49", 3.18, "28/09/2023 19:27", 4.43, "29/09/2023 15:47", "https://www.betexplorer.com/football/algeria/ligue-1/us-souf-oran/6qOsFaSf/")

$ws.Range("F17").Value = $row18[0]
$ws.Range("G17").Value = $row18[1]
$ws.Range("H17").Value = $row18[2]
$ws.Range("I17").Value = $row18[3]
$ws.Range("J17").Value = $row18[4]
$ws.Range("K17").Value = $row18[5]
$ws.Range("L17").Value = $row18[6]
$ws.Range("M17").Value = $row18[7]
$ws.Range("N17").Value = $row18[8]
$ws.Range("O17").Value = $row18[9]
$ws.Range("P17").Value = $row18[10]
$ws.Range("Q17").Value = $row18[11]
$ws.Range("R17").Value = $row18[12]
$ws.Range("S17").Value = $row18[13]
$ws.Range("T17").Value = $row18[14]
$ws.Range("U17").Value = $row18[15]
$ws.Range("V17").Value = $row18[16]

$ws.Range("F18").Value = $row17[0]
$ws.Range("G18").Value = $row17[1]
$ws.Range("H18").Value = $row17[2]
$ws.Range("I18").Value = $row17[3]
$ws.Range("J18").Value = $row17[4]
$ws.Range("K18").Value = $row17[5]
$ws.Range("L18").Value = $row17[6]
$ws.Range("M18").Value = $row17[7]
$ws.Range("N18").Value = $row17[8]
$ws.Range("O18").Value = $row17[9]
$ws.Range("P18").Value = $row17[10]
$ws.Range("Q18").Value = $row17[11]
$ws.Range("R18").Value = $row17[12]
$ws.Range("S18").Value = $row17[13]
$ws.Range("T18").Value = $row17[14]
$ws.Range("U18").Value = $row17[15]
$ws.Range("V18").Value = $row17[16]

# ---------------------------------------------------------------------------
# 3) Append four new match rows (31-34) at the bottom of the sheet, copying
#    the cell formatting pattern (bold/bordered index column, datetime format
#    on the match-date column) from the last existing data row (30).
# ---------------------------------------------------------------------------
$ws.Range("A30:V30").Copy()
$ws.Range("A31:V34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @(30, "algeria", "ligue-1", "2023-2024", 45240.64583333334, "Paradou", 1, "Oran", 0, 1.48, "19/10/2023 04:42", 1.46, "10/11/2023 14:57", 3.71, "19/10/2023 04:42", 3.91, "10/11/2023 14:57", 6.42, "19/10/2023 04:42", 8.55, "10/11/2023 14:57", "https://www.betexplorer.com/football/algeria/ligue-1/paradou-oran/ALKA1eA0/"),
    @(31, "algeria", "ligue-1", "2023-2024", 45240.64583333334, "US Souf", 3, "Constantine", 4, 3.14, "10/11/2023 06:42", 3.21, "10/11/2023 15:16", 2.82, "10/11/2023 06:42", 2.86, "10/11/2023 13:35", 2.5, "10/11/2023 06:42", 2.53, "10/11/2023 15:16", "https://www.betexplorer.com/football/algeria/ligue-1/us-souf-constantine/6mEJaZvD/"),
    @(32, "algeria", "ligue-1", "2023-2024", 45240.65625, "ES Setif", 1, "Kabylie", 0, 2.6, "09/11/2023 04:12", 2.06, "10/11/2023 12:22", 2.7, "09/11/2023 04:12", 3.06, "10/11/2023 13:50", 2.89, "09/11/2023 04:12", 3.93, "10/11/2023 12:22", "https://www.betexplorer.com/football/algeria/ligue-1/es-setif-kabylie/Cp59MdnQ/"),
    @(33, "algeria", "ligue-1", "2023-2024", 45240.70833333334, "USM Alger", 2, "CR Belouizdad", 1, 2.23, "09/11/2023 05:12", 2.35, "10/11/2023 16:57", 2.83, "09/11/2023 05:12", 2.82, "10/11/2023 16:57", 3.35, "09/11/2023 05:12", 3.6, "10/11/2023 16:57", "https://www.betexplorer.com/football/algeria/ligue-1/usm-alger-cr-belouizdad/IRsxIGfs/")
)

$targetRow = 31
foreach ($data in $newRows) {
    $ws.Cells.Item($targetRow, 1).Value = $data[0]
    $ws.Cells.Item($targetRow, 2).Value = $data[1]
    $ws.Cells.Item($targetRow, 3).Value = $data[2]
    $ws.Cells.Item($targetRow, 4).Value = $data[3]
    $ws.Cells.Item($targetRow, 5).Value = $data[4]
    $ws.Cells.Item($targetRow, 6).Value = $data[5]
    $ws.Cells.Item($targetRow, 7).Value = $data[6]
    $ws.Cells.Item($targetRow, 8).Value = $data[7]
    $ws.Cells.Item($targetRow, 9).Value = $data[8]
    $ws.Cells.Item($targetRow, 10).Value = $data[9]
    $ws.Cells.Item($targetRow, 11).Value = $data[10]
    $ws.Cells.Item($targetRow, 12).Value = $data[11]
    $ws.Cells.Item($targetRow, 13).Value = $data[12]
    $ws.Cells.Item($targetRow, 14).Value = $data[13]
    $ws.Cells.Item($targetRow, 15).Value = $data[14]
    $ws.Cells.Item($targetRow, 16).Value = $data[15]
    $ws.Cells.Item($targetRow, 17).Value = $data[16]
    $ws.Cells.Item($targetRow, 18).Value = $data[17]
    $ws.Cells.Item($targetRow, 19).Value = $data[18]
    $ws.Cells.Item($targetRow, 20).Value = $data[19]
    $ws.Cells.Item($targetRow, 21).Value = $data[20]
    $ws.Cells.Item($targetRow, 22).Value = $data[21]
    $targetRow = $targetRow + 1
}
